# topBot.xlsx edit: "made gui optional, added troop support"
# - Adds a new "bonus" column (AT) with the troop bonus for each
#   territory's continent.
# - Fixes the Africa continent group label (rows 22-27 had been
#   mis-tagged "Europe" in the AS/"group" column; now "Africa").
# - Extends conditional formatting to cover the new AT column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell AT1: "bonus" (copy header style from AS1, then set text) ---
$ws.Range("AS1").Copy() | Out-Null
$ws.Range("AT1").PasteSpecial(-4122) | Out-Null
$ws.Range("AT1").Value = "bonus"

# --- Fix mislabeled continent group for Africa rows (22-27): Europe -> Africa ---
$ws.Range("AS22:AS27").Value = "Africa"

# --- New AT column: troop bonus per continent, one value per territory row ---
$bonus = @{
    2  = 5;  3  = 5;  4  = 5;  5  = 5;  6  = 5;  7  = 5;  8  = 5;  9  = 5;  10 = 5;   # North America
    11 = 2;  12 = 2;  13 = 2;  14 = 2;                                               # South America
    15 = 5;  16 = 5;  17 = 5;  18 = 5;  19 = 5;  20 = 5;  21 = 5;                     # Europe
    22 = 3;  23 = 3;  24 = 3;  25 = 3;  26 = 3;  27 = 3;                             # Africa
    28 = 7;  29 = 7;  30 = 7;  31 = 7;  32 = 7;  33 = 7;  34 = 7;  35 = 7;  36 = 7;  37 = 7;  38 = 7;  39 = 7;  # Asia
    40 = 2;  41 = 2;  42 = 2;  43 = 2                                                # Australia
}

foreach ($row in 2..43) {
    $ws.Range("AT$row").Value = $bonus[$row]
}

# --- Extend conditional formatting range to include the new AT column ---
$ws.Range("B2:AR43 AT2:AT43").FormatConditions.Delete() | Out-Null

$cs = $ws.Range("B2:AR43 AT2:AT43").FormatConditions.AddColorScale(3)
$cs.ColorScaleCriteria.Item(1).Type = 0
$cs.ColorScaleCriteria.Item(1).Value = -1
$cs.ColorScaleCriteria.Item(1).FormatColor.Color = 255
$cs.ColorScaleCriteria.Item(2).Type = 0
$cs.ColorScaleCriteria.Item(2).Value = 0
$cs.ColorScaleCriteria.Item(2).FormatColor.Color = 20565
$cs.ColorScaleCriteria.Item(3).Type = 0
$cs.ColorScaleCriteria.Item(3).Value = 1
$cs.ColorScaleCriteria.Item(3).FormatColor.Color = 12611840
$cs.SetFirstPriority()

$rule2 = $ws.Range("B2:AR43 AT2:AT43").FormatConditions.Add(1, 3, "2")
$rule2.Interior.Color = 65535
$rule2.SetFirstPriority()

$rule3 = $ws.Range("B2:AR43 AT2:AT43").FormatConditions.Add(1, 3, "3")
$rule3.Interior.ColorIndex = -4142
$rule3.SetFirstPriority()

Write-Host "edit complete"
